$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value to G1 (new shared string "whoooo?")
$ws.Range("G1").Value = "whoooo?"

# Update the selection to G1, matching the diff's sheetView selection change
$ws.Range("G1").Select()
